# Add season-record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add labels for the three new columns ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the format from an existing header cell (e.g. A1) onto the
# newly added header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-49): every player on this roster shares the team's
# season record, so fill the same Wins/Losses/Ties values down the
# column for each row. ---
$ws.Range("AD2:AD49").Value = 76
$ws.Range("AE2:AE49").Value = 85
$ws.Range("AF2:AF49").Value = 0
